$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ANT1 Mid X / Mid Y values (row 2)
$ws.Range("B2").Value = 86.61
$ws.Range("C2").Value = -50.18

# Remove J2, J3, J4 rows (rows 20, 21, 22) - these components were dropped
$ws.Rows("20:22").Delete()
